$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 102, shifting existing rows 102:187 down to 103:188.
$ws.Rows("102:102").Insert()

# Populate the newly inserted row 102 with its data (Insert() does not
# duplicate the surrounding row's cell contents, only carries some format).
$ws.Range("A102").Value = 10
$ws.Range("B102").Value = "Vega Modelo de Temuco"
$ws.Range("C102").Value = "La Araucanía"
$ws.Range("D102").Value = 44587
$ws.Range("E102").Value = 9
$ws.Range("F102").Value = "Fruta"
$ws.Range("G102").Value = 100102
$ws.Range("H102").Value = "Cítricos"
$ws.Range("I102").Value = 100102006
$ws.Range("J102").Value = "Pomelo"
$ws.Range("K102").Value = "Start Ruby"
$ws.Range("L102").Value = "Primera"
$ws.Range("M102").Value = 80
$ws.Range("N102").Value = 14000
$ws.Range("O102").Value = 14000
$ws.Range("P102").Value = 14000
$ws.Range("Q102").Value = "$/bandeja 15 kilos granel"
$ws.Range("R102").Value = "Región de O'Higgins"
$ws.Range("S102").Value = 933
$ws.Range("T102").Value = 15
